$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.205287
$ws.Range("H2").Value = 6.615861000000001
$ws.Range("I2").Value = 0.08934890526417845
$ws.Range("J2").Value = 0.0954811485786439
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.009224
$ws.Range("N2").Value = 0.027672
$ws.Range("O2").Value = 0.01664709506025488
$ws.Range("P2").Value = 0.01697905168388481
$ws.Range("Q2").Value = 0.020341567288
$ws.Range("R2").Value = 0.183074105592
$ws.Range("S2").Value = 0.001487399719462486
$ws.Range("T2").Value = 0.001621179356553479
$ws.Range("G3").Value = 2.205287
$ws.Range("H3").Value = 6.615861000000001
$ws.Range("I3").Value = 0.08934890526417845
$ws.Range("J3").Value = 0.0954811485786439
$ws.Range("O3").Value = 0.9247000490894389
$ws.Range("P3").Value = 0.9431393206293146
$ws.Range("Q3").Value = 1.129917754520333
$ws.Range("R3").Value = 10.169259790683
$ws.Range("S3").Value = 0.08262093708387344
$ws.Range("T3").Value = 0.09005202560336885
$ws.Range("G4").Value = 2.205287
$ws.Range("H4").Value = 6.615861000000001
$ws.Range("I4").Value = 0.08934890526417845
$ws.Range("J4").Value = 0.0954811485786439
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.032499
$ws.Range("N4").Value = 0.064998
$ws.Range("O4").Value = 0.05865285585030608
$ws.Range("P4").Value = 0.03988162768680055
$ws.Range("Q4").Value = 0.07166962221300001
$ws.Range("R4").Value = 0.430017733278
$ws.Range("S4").Value = 0.005240568460842513
$ws.Range("T4").Value = 0.003807943618721562
$ws.Range("H5").Value = 50.93384900000001
$ws.Range("I5").Value = 0.687874737549802
$ws.Range("J5").Value = 0.7350853356881611
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.009224
$ws.Range("N5").Value = 0.027672
$ws.Range("O5").Value = 0.01664709506025488
$ws.Range("P5").Value = 0.01697905168388481
$ws.Range("Q5").Value = 0.1566046077253334
$ws.Range("R5").Value = 1.409441469528
$ws.Range("S5").Value = 0.01145111614553943
$ws.Range("T5").Value = 0.0124810519067151
$ws.Range("H6").Value = 50.93384900000001
$ws.Range("I6").Value = 0.687874737549802
$ws.Range("J6").Value = 0.7350853356881611
$ws.Range("O6").Value = 0.9247000490894389
$ws.Range("P6").Value = 0.9431393206293146
$ws.Range("Q6").Value = 8.698952455494114
$ws.Range("R6").Value = 78.29057209944702
$ws.Range("S6").Value = 0.6360778035796868
$ws.Range("T6").Value = 0.693287884105504
$ws.Range("H7").Value = 50.93384900000001
$ws.Range("I7").Value = 0.687874737549802
$ws.Range("J7").Value = 0.7350853356881611
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.032499
$ws.Range("N7").Value = 0.064998
$ws.Range("O7").Value = 0.05865285585030608
$ws.Range("P7").Value = 0.03988162768680055
$ws.Range("Q7").Value = 0.5517663862170001
$ws.Range("R7").Value = 3.310598317302001
$ws.Range("S7").Value = 0.04034581782457566
$ws.Range("T7").Value = 0.02931639967594205
$ws.Range("G8").Value = 0.32709
$ws.Range("H8").Value = 0.9812700000000001
$ws.Range("I8").Value = 0.01325230386015976
$ws.Range("J8").Value = 0.01416184328325004
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009224
$ws.Range("N8").Value = 0.027672
$ws.Range("O8").Value = 0.01664709506025488
$ws.Range("P8").Value = 0.01697905168388481
$ws.Range("Q8").Value = 0.00301707816
$ws.Range("R8").Value = 0.02715370344
$ws.Range("S8").Value = 0.0002206123621274622
$ws.Range("T8").Value = 0.0002404546690453793
$ws.Range("G9").Value = 0.32709
$ws.Range("H9").Value = 0.9812700000000001
$ws.Range("I9").Value = 0.01325230386015976
$ws.Range("J9").Value = 0.01416184328325004
$ws.Range("O9").Value = 0.9247000490894389
$ws.Range("P9").Value = 0.9431393206293146
$ws.Range("Q9").Value = 0.16759034009
$ws.Range("R9").Value = 1.50831306081
$ws.Range("S9").Value = 0.01225440603003789
$ws.Range("T9").Value = 0.01335659125302326
$ws.Range("G10").Value = 0.32709
$ws.Range("H10").Value = 0.9812700000000001
$ws.Range("I10").Value = 0.01325230386015976
$ws.Range("J10").Value = 0.01416184328325004
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.032499
$ws.Range("N10").Value = 0.064998
$ws.Range("O10").Value = 0.05865285585030608
$ws.Range("P10").Value = 0.03988162768680055
$ws.Range("Q10").Value = 0.01063009791
$ws.Range("R10").Value = 0.06378058746000001
$ws.Range("S10").Value = 0.0007772854679944051
$ws.Range("T10").Value = 0.0005647973611813953
$ws.Range("G11").Value = 4.7555295
$ws.Range("H11").Value = 9.511058999999999
$ws.Range("I11").Value = 0.1926739489130013
$ws.Range("J11").Value = 0.1372651023833856
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.009224
$ws.Range("N11").Value = 0.027672
$ws.Range("O11").Value = 0.01664709506025488
$ws.Range("P11").Value = 0.01697905168388481
$ws.Range("Q11").Value = 0.043865004108
$ws.Range("R11").Value = 0.263190024648
$ws.Range("S11").Value = 0.003207461543189424
$ws.Range("T11").Value = 0.002330631267761244
$ws.Range("G12").Value = 4.7555295
$ws.Range("H12").Value = 9.511058999999999
$ws.Range("I12").Value = 0.1926739489130013
$ws.Range("J12").Value = 0.1372651023833856
$ws.Range("O12").Value = 0.9247000490894389
$ws.Range("P12").Value = 0.9431393206293146
$ws.Range("Q12").Value = 2.4365795536795
$ws.Range("R12").Value = 14.619477322077
$ws.Range("S12").Value = 0.1781656100181083
$ws.Range("T12").Value = 0.1294601154079796
$ws.Range("G13").Value = 4.7555295
$ws.Range("H13").Value = 9.511058999999999
$ws.Range("I13").Value = 0.1926739489130013
$ws.Range("J13").Value = 0.1372651023833856
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.5
$ws.Range("M13").Value = 0.032499
$ws.Range("N13").Value = 0.064998
$ws.Range("O13").Value = 0.05865285585030608
$ws.Range("P13").Value = 0.03988162768680055
$ws.Range("Q13").Value = 0.1545499532205
$ws.Range("R13").Value = 0.618199812882
$ws.Range("S13").Value = 0.0113008773517035
$ws.Range("T13").Value = 0.005474355707644744
$ws.Range("G14").Value = 0.4158899999999999
$ws.Range("H14").Value = 1.24767
$ws.Range("I14").Value = 0.01685010441285836
$ws.Range("J14").Value = 0.01800657006655923
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.009224
$ws.Range("N14").Value = 0.027672
$ws.Range("O14").Value = 0.01664709506025488
$ws.Range("P14").Value = 0.01697905168388481
$ws.Range("Q14").Value = 0.003836169359999999
$ws.Range("R14").Value = 0.03452552423999999
$ws.Range("S14").Value = 0.0002805052899360733
$ws.Range("T14").Value = 0.0003057344838096022
$ws.Range("G15").Value = 0.4158899999999999
$ws.Range("H15").Value = 1.24767
$ws.Range("I15").Value = 0.01685010441285836
$ws.Range("J15").Value = 0.01800657006655923
$ws.Range("O15").Value = 0.9247000490894389
$ws.Range("P15").Value = 0.9431393206293146
$ws.Range("Q15").Value = 0.21308858889
$ws.Range("R15").Value = 1.91779730001
$ws.Range("S15").Value = 0.01558129237773229
$ws.Range("T15").Value = 0.01698270425943882
$ws.Range("G16").Value = 0.4158899999999999
$ws.Range("H16").Value = 1.24767
$ws.Range("I16").Value = 0.01685010441285836
$ws.Range("J16").Value = 0.01800657006655923
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.5
$ws.Range("M16").Value = 0.032499
$ws.Range("N16").Value = 0.064998
$ws.Range("O16").Value = 0.05865285585030608
$ws.Range("P16").Value = 0.03988162768680055
$ws.Range("Q16").Value = 0.01351600911
$ws.Range("R16").Value = 0.08109605465999999
$ws.Range("S16").Value = 0.0009883067451899876
$ws.Range("T16").Value = 0.0007181313233108025
